$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Paragraph "Model2:" -> "Using transfer learning from " + "Model2" ---
$r1 = $d.Content
[void]$r1.Find.Execute("Model2:")
$p1 = $r1.Paragraphs(1)
$xml1 = '<w:p ' + $wNs + '><w:r><w:t xml:space="preserve">Using transfer learning from </w:t></w:r><w:r><w:t>Model2</w:t></w:r></w:p>'
[void]$p1.Range.InsertXML($xml1)

# --- Paragraph "0 [en dash] 99.96" -> Mse block + extra results ---
$r2 = $d.Content
[void]$r2.Find.Execute([string][char]0x0030 + " " + [string][char]0x2013 + " 99.96")
$p2 = $r2.Paragraphs(1)

$xml2 = (
  '<w:p ' + $wNs + '><w:r><w:t>Mse:</w:t></w:r></w:p>' +
  '<w:p ' + $wNs + '><w:r><w:t xml:space="preserve">0 </w:t></w:r><w:r><w:t>-</w:t></w:r><w:r><w:t xml:space="preserve"> 99.96</w:t></w:r></w:p>' +
  '<w:p ' + $wNs + '><w:r><w:t xml:space="preserve">1 </w:t></w:r><w:r><w:t>-</w:t></w:r><w:r><w:t xml:space="preserve"> 104</w:t></w:r></w:p>' +
  '<w:p ' + $wNs + '><w:r><w:t xml:space="preserve">2 </w:t></w:r><w:r><w:t>-</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>100</w:t></w:r></w:p>' +
  '<w:p ' + $wNs + '><w:r><w:t>3 - 98</w:t></w:r></w:p>' +
  '<w:p ' + $wNs + '><w:r><w:t>5</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>-</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>96</w:t></w:r></w:p>' +
  '<w:p ' + $wNs + '><w:r><w:t>vs</w:t></w:r></w:p>' +
  '<w:p ' + $wNs + '><w:r><w:t>0 to 5 - 93</w:t></w:r><w:r><w:t xml:space="preserve"> - saved as model4v1</w:t></w:r></w:p>' +
  '<w:p ' + $wNs + '><w:r><w:t>5 to 10 - 90 - saved as model4v2</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
)
[void]$p2.Range.InsertXML($xml2)
